$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ma" -> "ma🔑" for the primary key columns in each mini-table
$ws.Range("A5").Value = "ma🔑"
$ws.Range("C5").Value = "ma🔑"
$ws.Range("E5").Value = "ma🔑"
$ws.Range("G5").Value = "ma🔑"
$ws.Range("A15").Value = "ma🔑"

# "ma_cap_do" -> "ma_cap_do🔑"
$ws.Range("E16").Value = "ma_cap_do🔑"

# Update the selected cell/range shown when the sheet is opened
$ws.Range("K16").Select() | Out-Null
